$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.319081783294678
$ws.Range("B1").Value = 3.692748546600342
$ws.Range("C1").Value = 3.097846746444702
$ws.Range("D1").Value = 2.03364109992981
$ws.Range("E1").Value = 1.169101357460022
